$wb = $excel.ActiveWorkbook

# Rename the "Example" sheet to "Example.In.Excel" to match the renamed
# Excel document / script name (per commit message).
$wsExample = $wb.Worksheets.Item("Example")
$wsExample.Name = "Example.In.Excel"

# Update the selection/scroll state on "Observed.Matrix" first so that the
# final Activate() below (on "Example.In.Excel") is what leaves the
# tabSelected flag on the Example sheet, matching the target workbook.
$wsObserved = $wb.Worksheets.Item("Observed.Matrix")
$wsObserved.Activate()
$wsObserved.Range("E18").Select()

# Update the selection/scroll state on the renamed Example sheet.
$wsExample.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$wsExample.Range("B36").Select()
